# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.333.35"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.565.02"
$ws.Range("E3").Value = "  +0.00%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "289.57"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3758"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.30%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "49.14"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.3361"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07504"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.07%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.125"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.82%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.36%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "20.78"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.09%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.894"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.59%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.856"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").Value = "1.564.90"
$ws.Range("E16").Value = "  +0.03%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001114"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.50%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "89.23"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.00%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06705"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.53%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.45%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.160"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.43%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "16.32"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").Value = "22.324.02"
$ws.Range("E24").Value = "  -0.19%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.376"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.52%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.657"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -4.94%  "
$ws.Range("E27").Value = "  -0.77%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "147.39"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.38%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "4.989"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.20%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "125.09"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").Value = "1.736.50"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.016"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.16%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.9788"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.15%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.930"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.48%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "9.826"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.08436"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.404"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +9.67%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02447"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.46%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.2260"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.80%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.06384"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.69%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "5.333"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.68%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.6225"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.83%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "10.94"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -6.59%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.61%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.85"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.87%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.790"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.75%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.5801"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.96%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.042"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.63%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.246"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.13%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "124.10"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.41%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.07306"
$c.Style = "Normal"
